# RoutingPlan.xlsx update
#
# 1. Insert a blank spacer row after "/register" (POST) block, before "/admin" block.
# 2. Insert a blank spacer row after the "/admin/user/:id/profile/edit" (POST) block,
#    before the "/user/:id/course" block.
# 3. Mark the "routers" (col F) progress status for each route:
#      - rows 2..5  (landing/login/register)          -> "v"    (done)
#      - rows 6..14 (admin / admin course routes)      -> "U"    (in progress)
#      - rows 15,16 (admin/user profile edit routes)   -> "XXX" + "NANTI" note (col H)
#      - rows 22,23 (user profile edit routes)         -> "XXX" + "NANTI" note (col H)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert spacer row before the "/admin" row (currently row 9) ---
$ws.Range("B9").EntireRow.Insert()
$ws.Range("B8").Copy()
$ws.Range("B9").PasteSpecial(-4122)   # xlPasteFormats

# --- Insert spacer row before the "/user/:id/course" row (currently row 21 after the shift above) ---
$ws.Range("B21").EntireRow.Insert()
$ws.Range("B20").Copy()
$ws.Range("B21").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0

# --- Column F progress markers ---
# NB: new shared-string entries are created in the order cells are written
# (NANTI, then U, then XXX) so the shared-string table matches the target.
$ws.Range("H19").Value2 = "NANTI"
$ws.Range("H20").Value2 = "NANTI"
$ws.Range("H27").Value2 = "NANTI"
$ws.Range("H28").Value2 = "NANTI"

$ws.Range("F10").Value2 = "U"
$ws.Range("F11").Value2 = "U"
$ws.Range("F12").Value2 = "U"
$ws.Range("F13").Value2 = "U"
$ws.Range("F14").Value2 = "U"
$ws.Range("F15").Value2 = "U"
$ws.Range("F16").Value2 = "U"
$ws.Range("F17").Value2 = "U"
$ws.Range("F18").Value2 = "U"

$ws.Range("F19").Value2 = "XXX"
$ws.Range("F20").Value2 = "XXX"
$ws.Range("F27").Value2 = "XXX"
$ws.Range("F28").Value2 = "XXX"

$ws.Range("F5").Value2 = "v"
$ws.Range("F6").Value2 = "v"
$ws.Range("F7").Value2 = "v"
$ws.Range("F8").Value2 = "v"

# --- Column E width ---
$ws.Columns.Item(5).ColumnWidth = 54.42578125

# --- Selection ---
$ws.Range("E18").Select()
